$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as Text so numeric-looking strings
# (e.g. "0.9277", "274.05") are kept as literal text instead of
# being auto-converted to numbers -- matches the source data which
# stores prices as inline text.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Coin name / link swap (rows 40-41) ---
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "20.250.12"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.443.85"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  +0.81%  "
$ws.Range("D5").Value = "0.9277"
$ws.Range("E5").Value = "  -7.21%  "
$ws.Range("D6").Value = "274.05"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").Value = "0.3634"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "39.48"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "1.024"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "0.06514"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "0.9983"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "5.348"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").Value = "17.49"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "6.057"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "0.00001010"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "1.440.16"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "0.9471"
$ws.Range("E18").Value = "  -5.19%  "
$ws.Range("D19").Value = "0.05657"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "68.78"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "5.376"
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("D22").Value = "14.23"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "10.76"
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("D24").Value = "2.245"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "20.270.28"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").Value = "140.12"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("D27").Value = "2.040"
$ws.Range("E27").Value = "  -10.11%  "
$ws.Range("D28").Value = "16.94"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").Value = "1.593.25"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").Value = "110.60"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "3.990"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").Value = "4.803"
$ws.Range("E32").Value = "  -10.54%  "
$ws.Range("D33").Value = "0.7835"
$ws.Range("E33").Value = "  -5.52%  "
$ws.Range("D34").Value = "0.07685"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "1.452"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").Value = "0.05653"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("D37").Value = "4.648"
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").Value = "1.117"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").Value = "0.01992"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("D40").Value = "0.9383"
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("D41").Value = "10.17"
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").Value = "0.1840"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("D43").Value = "6.926"
$ws.Range("E43").Value = "  -18.34%  "
$ws.Range("D44").Value = "0.5201"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "3.469"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").Value = "11.73"
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("D47").Value = "115.46"
$ws.Range("D48").Value = "0.5103"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "1.725"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").Value = "0.06365"
$ws.Range("D51").Value = "0.9833"
$ws.Range("E51").Value = "  -1.62%  "

# Restore the Price column to its original (General/default) style
# now that the text values are safely stored, so no stray
# number-format style survives in the saved workbook.
$ws.Range("D2:D51").ClearFormats()

